$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M: company names, each linked back to the company's website
# (mirrors the existing "email" hyperlink column C) -- lets the import
# test read multiple related values (name + link) out of one column.
$ws.Range("M1").Value = "Maatwebsite"
$ws.Range("M2").Value = "Laravel"

$ws.Hyperlinks.Add($ws.Range("M1"), "https://www.maatwebsite.nl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("M2"), "https://laravel.com") | Out-Null

# Hyperlinks.Add reapplies its own flavour of the "Hyperlink" cell style;
# put the cells back on the same shared style C1/C2 already use.
$ws.Range("M1").Style = $ws.Range("C1").Style
$ws.Range("M2").Style = $ws.Range("C2").Style

# Give the new column a sensible width, like the other text columns.
$ws.Columns.Item(13).ColumnWidth = 10

# Restore the (arbitrary) last-used selection.
$ws.Range("F10").Select() | Out-Null
